$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = "stimuli/img_anzgh.png"
$ws.Range("M2").Value = 75.10526315789474
$ws.Range("N2").Value = 55.76315789473684
$ws.Range("O2").Value = 65.43421052631579
$ws.Range("P2").Value = 38
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 6

# Row 3
$ws.Range("L3").Value = "stimuli/img_cmyvx.png"
$ws.Range("M3").Value = 64.25
$ws.Range("N3").Value = 40.09375
$ws.Range("O3").Value = 52.171875
$ws.Range("P3").Value = 32
$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 4
$ws.Range("S3").Value = 4

# Row 4
$ws.Range("L4").Value = "stimuli/img_2pnl2.png"
$ws.Range("M4").Value = 6.621621621621622
$ws.Range("N4").Value = 7.135135135135135
$ws.Range("O4").Value = 6.878378378378379
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1

# Row 5
$ws.Range("L5").Value = "stimuli/img_aweye.png"
$ws.Range("M5").Value = 53.42105263157895
$ws.Range("N5").Value = 31.84210526315789
$ws.Range("O5").Value = 42.63157894736842
$ws.Range("P5").Value = 38

# Row 6
$ws.Range("L6").Value = "stimuli/img_gbypq.png"
$ws.Range("M6").Value = 76.27500000000001
$ws.Range("N6").Value = 51.925
$ws.Range("O6").Value = 64.09999999999999
$ws.Range("P6").Value = 40
$ws.Range("Q6").Value = 6
$ws.Range("R6").Value = 6
$ws.Range("S6").Value = 6

# Row 7
$ws.Range("L7").Value = "stimuli/img_t4hvr.png"
$ws.Range("M7").Value = 61.69230769230769
$ws.Range("N7").Value = 39.76923076923077
$ws.Range("O7").Value = 50.73076923076923
$ws.Range("P7").Value = 39
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 3
$ws.Range("S7").Value = 3

# Row 8
$ws.Range("L8").Value = "stimuli/img_3bxjb.png"
$ws.Range("M8").Value = 87.28571428571429
$ws.Range("N8").Value = 72.65714285714286
$ws.Range("O8").Value = 79.97142857142858
$ws.Range("P8").Value = 35
$ws.Range("Q8").Value = 10
$ws.Range("R8").Value = 10
$ws.Range("S8").Value = 10

# Row 9
$ws.Range("H9").Value = "bedrooms"
$ws.Range("I9").Value = "target"
$ws.Range("K9").Value = "j"
$ws.Range("L9").Value = "stimuli/img_ic3os.png"
$ws.Range("M9").Value = 84.79069767441861
$ws.Range("N9").Value = 66.16279069767442
$ws.Range("O9").Value = 75.47674418604652
$ws.Range("P9").Value = 43
$ws.Range("Q9").Value = 9
$ws.Range("R9").Value = 9
$ws.Range("S9").Value = 9

# Row 11
$ws.Range("H11").Value = "bedrooms"
$ws.Range("I11").Value = "target"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_juob3.png"
$ws.Range("M11").Value = 79.92105263157895
$ws.Range("N11").Value = 59.78947368421053
$ws.Range("O11").Value = 69.85526315789474
$ws.Range("P11").Value = 38
$ws.Range("Q11").Value = 7
$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 7

# Row 12
$ws.Range("H12").Value = "living_rooms"
$ws.Range("L12").Value = "stimuli/img_pbsj1.png"
$ws.Range("M12").Value = 73.88636363636364
$ws.Range("N12").Value = 51.52272727272727
$ws.Range("O12").Value = 62.70454545454545
$ws.Range("P12").Value = 44
$ws.Range("Q12").Value = 6
$ws.Range("R12").Value = 6
$ws.Range("S12").Value = 6

# Row 13
$ws.Range("L13").Value = "stimuli/img_9pfbj.png"
$ws.Range("M13").Value = 91.27272727272727
$ws.Range("N13").Value = 80.09090909090909
$ws.Range("O13").Value = 85.68181818181819
$ws.Range("P13").Value = 33
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 10

# Row 14
$ws.Range("L14").Value = "stimuli/img_72fmj.png"
$ws.Range("M14").Value = 53.87179487179487
$ws.Range("N14").Value = 36.02564102564103
$ws.Range("O14").Value = 44.94871794871795
$ws.Range("P14").Value = 39
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 3
$ws.Range("S14").Value = 3

# Row 15
$ws.Range("H15").Value = "bedrooms"
$ws.Range("I15").Value = "target"
$ws.Range("K15").Value = "j"
$ws.Range("L15").Value = "stimuli/img_jivhq.png"
$ws.Range("M15").Value = 37
$ws.Range("N15").Value = 22.26530612244898
$ws.Range("O15").Value = 29.63265306122449
$ws.Range("P15").Value = 49
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = 2
$ws.Range("S15").Value = 2

# Row 16
$ws.Range("H16").Value = "living_rooms"
$ws.Range("I16").Value = "distractor"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_zxvl3.png"
$ws.Range("M16").Value = 68.78260869565217
$ws.Range("N16").Value = 47.56521739130435
$ws.Range("O16").Value = 58.17391304347827
$ws.Range("P16").Value = 46
$ws.Range("Q16").Value = 5
$ws.Range("R16").Value = 5
$ws.Range("S16").Value = 5

# Row 17
$ws.Range("L17").Value = "stimuli/img_yteqw.png"
$ws.Range("M17").Value = 66.83783783783784
$ws.Range("N17").Value = 43.78378378378378
$ws.Range("O17").Value = 55.31081081081081
$ws.Range("P17").Value = 37

# Row 18
$ws.Range("L18").Value = "stimuli/img_f4jxo.png"
$ws.Range("M18").Value = 82.91666666666667
$ws.Range("N18").Value = 65.52777777777777
$ws.Range("O18").Value = 74.22222222222223
$ws.Range("R18").Value = 8
$ws.Range("S18").Value = 8

# Row 19
$ws.Range("L19").Value = "stimuli/img_kzg3h.png"
$ws.Range("M19").Value = 77.02777777777777
$ws.Range("N19").Value = 56.22222222222222
$ws.Range("O19").Value = 66.625
$ws.Range("P19").Value = 36
$ws.Range("Q19").Value = 7
$ws.Range("R19").Value = 7
$ws.Range("S19").Value = 7

# Row 20
$ws.Range("H20").Value = "kitchens"
$ws.Range("I20").Value = "distractor"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_eppte.png"
$ws.Range("M20").Value = 78.42424242424242
$ws.Range("N20").Value = 57.03030303030303
$ws.Range("O20").Value = 67.72727272727272
$ws.Range("P20").Value = 33
$ws.Range("Q20").Value = 7
$ws.Range("R20").Value = 7
$ws.Range("S20").Value = 7

# Row 21
$ws.Range("L21").Value = "stimuli/img_z3yzz.png"
$ws.Range("M21").Value = 71.71052631578948
$ws.Range("N21").Value = 49.81578947368421
$ws.Range("O21").Value = 60.76315789473685
$ws.Range("P21").Value = 38
$ws.Range("Q21").Value = 5
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 5

# Row 22
$ws.Range("L22").Value = "stimuli/img_ose78.png"
$ws.Range("M22").Value = 80.19444444444444
$ws.Range("N22").Value = 60.25
$ws.Range("O22").Value = 70.22222222222223
$ws.Range("Q22").Value = 8

# Row 23
$ws.Range("H23").Value = "kitchens"
$ws.Range("I23").Value = "distractor"
$ws.Range("K23").Value = "f"
$ws.Range("L23").Value = "stimuli/img_ps986.png"
$ws.Range("M23").Value = 90.46428571428571
$ws.Range("N23").Value = 70.60714285714286
$ws.Range("O23").Value = 80.53571428571428
$ws.Range("P23").Value = 28
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = 10
$ws.Range("S23").Value = 10

# Row 24
$ws.Range("H24").Value = "living_rooms"
$ws.Range("I24").Value = "distractor"
$ws.Range("K24").Value = "f"
$ws.Range("L24").Value = "stimuli/img_95hiv.png"
$ws.Range("M24").Value = 84.04545454545455
$ws.Range("N24").Value = 67.31818181818181
$ws.Range("O24").Value = 75.68181818181819
$ws.Range("P24").Value = 44
$ws.Range("Q24").Value = 9
$ws.Range("R24").Value = 9
$ws.Range("S24").Value = 9

# Row 25
$ws.Range("L25").Value = "stimuli/img_fnu4h.png"
$ws.Range("M25").Value = 85.87179487179488
$ws.Range("N25").Value = 70.71794871794872
$ws.Range("O25").Value = 78.2948717948718
$ws.Range("P25").Value = 39
$ws.Range("Q25").Value = 9
$ws.Range("R25").Value = 9
$ws.Range("S25").Value = 9

# Row 26
$ws.Range("H26").Value = "bedrooms"
$ws.Range("I26").Value = "target"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_cgdyc.png"
$ws.Range("M26").Value = 32.93023255813954
$ws.Range("N26").Value = 14.04651162790698
$ws.Range("O26").Value = 23.48837209302326
$ws.Range("P26").Value = 43
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 1

# Row 27
$ws.Range("L27").Value = "stimuli/img_1vq1v.png"
$ws.Range("M27").Value = 69.42857142857143
$ws.Range("N27").Value = 46.59523809523809
$ws.Range("O27").Value = 58.01190476190476
$ws.Range("P27").Value = 42
$ws.Range("Q27").Value = 5
$ws.Range("R27").Value = 5
$ws.Range("S27").Value = 5
